$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'20.453.19"
$ws.Range("E2").Value2 = "  -7.20%  "

$ws.Range("D3").Value2 = "'1.446.19"
$ws.Range("E3").Value2 = "  -6.99%  "

$ws.Range("D4").Value2 = "'1.005"
$ws.Range("E4").Value2 = "  +0.42%  "

$ws.Range("D5").Value2 = "'1.007"
$ws.Range("E5").Value2 = "  +0.51%  "

$ws.Range("D6").Value2 = "'278.11"
$ws.Range("E6").Value2 = "  -4.37%  "

$ws.Range("D7").Value2 = "'0.3719"
$ws.Range("E7").Value2 = "  -5.29%  "

$ws.Range("D8").Value2 = "'0.3068"
$ws.Range("E8").Value2 = "  -4.66%  "

$ws.Range("D9").Value2 = "'40.88"
$ws.Range("E9").Value2 = "  -7.80%  "

$ws.Range("D10").Value2 = "'1.012"
$ws.Range("E10").Value2 = "  -5.88%  "

$ws.Range("D11").Value2 = "'0.06536"
$ws.Range("E11").Value2 = "  -9.31%  "

$ws.Range("D12").Value2 = "'1.009"
$ws.Range("E12").Value2 = "  +0.77%  "

$ws.Range("D13").Value2 = "'5.388"
$ws.Range("E13").Value2 = "  -4.96%  "

$ws.Range("D14").Value2 = "'17.16"
$ws.Range("E14").Value2 = "  -8.31%  "

$ws.Range("B15").Value2 = "Chainlink"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value2 = "'6.130"
$ws.Range("E15").Value2 = "  -8.45%  "

$ws.Range("B16").Value2 = "WrappedEther"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value2 = "'1.448.67"
$ws.Range("E16").Value2 = "  -6.64%  "

$ws.Range("D17").Value2 = "'0.00001015"
$ws.Range("E17").Value2 = "  -7.82%  "

$ws.Range("D18").Value2 = "'76.52"
$ws.Range("E18").Value2 = "  -8.32%  "

$ws.Range("D19").Value2 = "'0.05861"
$ws.Range("E19").Value2 = "  -11.02%  "

$ws.Range("D21").Value2 = "'5.713"
$ws.Range("E21").Value2 = "  -8.25%  "

$ws.Range("D22").Value2 = "'14.37"
$ws.Range("E22").Value2 = "  -7.14%  "

$ws.Range("D23").Value2 = "'10.89"
$ws.Range("E23").Value2 = "  -2.90%  "

$ws.Range("D24").Value2 = "'2.292"
$ws.Range("E24").Value2 = "  -2.91%  "

$ws.Range("D25").Value2 = "'20.436.99"
$ws.Range("E25").Value2 = "  -7.31%  "

$ws.Range("D26").Value2 = "'143.19"
$ws.Range("E26").Value2 = "  -3.25%  "

$ws.Range("D27").Value2 = "'2.220"
$ws.Range("E27").Value2 = "  -7.03%  "

$ws.Range("D28").Value2 = "'17.00"
$ws.Range("E28").Value2 = "  -8.21%  "

$ws.Range("D29").Value2 = "'1.612.92"
$ws.Range("E29").Value2 = "  -6.64%  "

$ws.Range("D30").Value2 = "'109.08"
$ws.Range("E30").Value2 = "  -8.01%  "

$ws.Range("D31").Value2 = "'0.9109"
$ws.Range("E31").Value2 = "  -7.79%  "

$ws.Range("B32").Value2 = "HuobiToken"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value2 = "'3.651"
$ws.Range("E32").Value2 = "  -25.11%  "

$ws.Range("B33").Value2 = "Filecoin"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value2 = "'5.390"
$ws.Range("E33").Value2 = "  -8.41%  "

$ws.Range("D34").Value2 = "'0.07731"
$ws.Range("E34").Value2 = "  -7.00%  "

$ws.Range("D35").Value2 = "'8.293"
$ws.Range("E35").Value2 = "  -9.53%  "

$ws.Range("D36").Value2 = "'1.006"
$ws.Range("E36").Value2 = "  +0.57%  "

$ws.Range("B37").Value2 = "WEMIXTOKEN"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value2 = "'1.411"
$ws.Range("E37").Value2 = "  -12.48%  "

$ws.Range("B38").Value2 = "Aptos"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value2 = "'10.72"
$ws.Range("E38").Value2 = "  -0.22%  "

$ws.Range("D39").Value2 = "'0.05594"
$ws.Range("E39").Value2 = "  -7.09%  "

$ws.Range("D40").Value2 = "'4.726"
$ws.Range("E40").Value2 = "  -7.48%  "

$ws.Range("E41").Value2 = "  -6.84%  "

$ws.Range("D42").Value2 = "'0.02039"
$ws.Range("E42").Value2 = "  -9.89%  "

$ws.Range("D43").Value2 = "'0.1912"
$ws.Range("E43").Value2 = "  -7.04%  "

$ws.Range("D44").Value2 = "'3.585"
$ws.Range("E44").Value2 = "  -4.51%  "

$ws.Range("D45").Value2 = "'0.5314"
$ws.Range("E45").Value2 = "  -8.53%  "

$ws.Range("D46").Value2 = "'11.99"
$ws.Range("E46").Value2 = "  -8.23%  "

$ws.Range("D47").Value2 = "'0.5146"
$ws.Range("E47").Value2 = "  -7.67%  "

$ws.Range("D48").Value2 = "'110.95"
$ws.Range("E48").Value2 = "  -5.53%  "

$ws.Range("D49").Value2 = "'1.773"
$ws.Range("E49").Value2 = "  -5.77%  "

$ws.Range("D50").Value2 = "'1.054"
$ws.Range("E50").Value2 = "  -7.16%  "

$ws.Range("D51").Value2 = "'1.008"
$ws.Range("E51").Value2 = "  +0.61%  "
